$wb = $excel.ActiveWorkbook
$hub = $wb.Worksheets.Item("hub")

# Insert two new columns (H:I) on the "hub" sheet, pushing the old H..M to J..O.
$hub.Columns("H:I").Insert() | Out-Null
$hub.Columns("H:I").ColumnWidth = 11.33

# New header cells for the inserted columns.
$hub.Range("H1").Value = "width"
$hub.Range("I1").Value = "hight"

# Fill the new "width"/"hight" columns with "70px" for the two data rows that
# describe a sizeable component (rows 5 and 6).
$hub.Range("H5").Value = "70px"
$hub.Range("I5").Value = "70px"
$hub.Range("H6").Value = "70px"
$hub.Range("I6").Value = "70px"

# Drop the stale "testing123" values that used to live in column B of rows 5/6.
$hub.Range("B5").Value = ""
$hub.Range("B6").Value = ""

# Rows 5/6 also lose their "televisie" (name_en) entry -- after the column
# insert above, that value now sits in column L.
$hub.Range("L5").Value = ""
$hub.Range("L6").Value = ""

# Make "hub" the active tab/sheet with the given selection...
$hub.Activate()
$hub.Range("I25").Select() | Out-Null

# ...and make sure "system" is no longer flagged as the selected tab (it keeps
# its own prior cell selection, it just stops being the active sheet).
$system = $wb.Worksheets.Item("system")
$system.Range("A3").Select() | Out-Null
$hub.Activate()
